$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "Q8" in J1, copying style/format from existing header I1
$ws.Range("I1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null
$ws.Range("J1").Value = "Q8"

# Update simulated/bugfixed numeric data for existing and new cells
$ws.Range("B2").Value = 0.2298834519089068
$ws.Range("C2").Value = 0.3255607007654389
$ws.Range("D2").Value = 0.06513869270445388
$ws.Range("E2").Value = 0.9005966785379966
$ws.Range("F2").Value = 0.8005739094407884
$ws.Range("G2").Value = 0.4467158494068597
$ws.Range("H2").Value = 0.5949089188422247
$ws.Range("B3").Value = 0.2402751528550954
$ws.Range("C3").Value = -0.0201468552058896
$ws.Range("D3").Value = 0.815311130627653
$ws.Range("E3").Value = 0.7152883615304448
$ws.Range("F3").Value = 0.3614303014965162
$ws.Range("G3").Value = 0.5096233709318814
$ws.Range("B4").Value = -0.3601468552058896
$ws.Range("C4").Value = 0.4753111306276531
$ws.Range("D4").Value = 0.3752883615304449
$ws.Range("E4").Value = 0.02143030149651622
$ws.Range("F4").Value = 0.1696233709318813
$ws.Range("G4").Value = 0.3270570324401456
$ws.Range("H4").Value = -0.1735990538361767
$ws.Range("I4").Value = 0.10675723343443
$ws.Range("J4").Value = -0.06520887812495521
$ws.Range("B5").Value = 0.1453111306276531
$ws.Range("C5").Value = 0.04528836153044491
$ws.Range("D5").Value = -0.3085696985034838
$ws.Range("E5").Value = -0.1603766290681187
$ws.Range("F5").Value = -0.002942967559854404
$ws.Range("G5").Value = -0.5035990538361766
$ws.Range("H5").Value = -0.22324276656557
$ws.Range("I5").Value = -0.3952088781249552
$ws.Range("B6").Value = 0.1274104041151531
$ws.Range("C6").Value = -0.2264476559187756
$ws.Range("D6").Value = -0.07825458648341049
$ws.Range("E6").Value = 0.07917907502485377
$ws.Range("F6").Value = -0.4214770112514685
$ws.Range("G6").Value = -0.1411207239808618
$ws.Range("H6").Value = -0.313086835540247
$ws.Range("B7").Value = 0.1020781346750965
$ws.Range("C7").Value = 0.2502712041104616
$ws.Range("D7").Value = 0.4077048656187259
$ws.Range("E7").Value = -0.0929512206575964
$ws.Range("F7").Value = 0.1874050666130103
$ws.Range("G7").Value = 0.01543895505362509
$ws.Range("B8").Value = -0.05037662906811868
$ws.Range("C8").Value = 0.1070570324401456
$ws.Range("D8").Value = -0.3935990538361767
$ws.Range("E8").Value = -0.11324276656557
$ws.Range("F8").Value = -0.2852088781249552
$ws.Range("G8").Value = 0.02053824020493492
$ws.Range("H8").Value = -0.5916519254275203
$ws.Range("I8").Value = 0.09668430691196023
$ws.Range("B9").Value = 0.3310570324401456
$ws.Range("C9").Value = -0.1695990538361767
$ws.Range("D9").Value = 0.11075723343443
$ws.Range("E9").Value = -0.06120887812495521
$ws.Range("F9").Value = 0.2445382402049349
$ws.Range("G9").Value = -0.3676519254275204
$ws.Range("H9").Value = 0.3206843069119602
$ws.Range("B10").Value = -0.2961090942996805
$ws.Range("C10").Value = -0.01575280702907378
$ws.Range("D10").Value = -0.187718918588459
$ws.Range("E10").Value = 0.1180281997414312
$ws.Range("F10").Value = -0.4941619658910241
$ws.Range("G10").Value = 0.1941742664484565
$ws.Range("B11").Value = 0.3250097553919601
$ws.Range("C11").Value = 0.153043643832575
$ws.Range("D11").Value = 0.4587907621624651
$ws.Range("E11").Value = -0.1533994034699902
$ws.Range("F11").Value = 0.5349368288694905
$ws.Range("B12").Value = -0.1523844033498989
$ws.Range("C12").Value = 0.1533627149799912
$ws.Range("D12").Value = -0.4588274506524641
$ws.Range("E12").Value = 0.2295087816870165
$ws.Range("B13").Value = 0.4147680207538826
$ws.Range("C13").Value = -0.1974221448785727
$ws.Range("D13").Value = 0.4909140874609079
$ws.Range("B14").Value = -0.1722147680129069
$ws.Range("C14").Value = 0.5161214643265737
$ws.Range("B15").Value = 0.3238631410950035
